$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '62.803.30'
$ws.Range('E2').Value = '  -8.04%  '
$ws.Range('D3').Value = '3.215.48'
$ws.Range('E3').Value = '  -10.65%  '
$ws.Range('E4').Value = '  +0.02%  '
$ws.Range('D5').Value = "'172.58"
$ws.Range('D5').ClearFormats()
$ws.Range('E5').Value = '  -16.87%  '
$ws.Range('D6').Value = "'505.55"
$ws.Range('D6').ClearFormats()
$ws.Range('E6').Value = '  -11.32%  '
$ws.Range('D7').Value = "'0.582"
$ws.Range('D7').ClearFormats()
$ws.Range('E7').Value = '  -5.09%  '
$ws.Range('E8').Value = '  +0.09%  '
$ws.Range('D9').Value = '3.207.46'
$ws.Range('E9').Value = '  -10.67%  '
$ws.Range('D10').Value = "'0.600"
$ws.Range('D10').ClearFormats()
$ws.Range('E10').Value = '  -12.41%  '
$ws.Range('D11').Value = "'55.53"
$ws.Range('D11').ClearFormats()
$ws.Range('E11').Value = '  -13.19%  '
$ws.Range('D12').Value = "'0.127"
$ws.Range('D12').ClearFormats()
$ws.Range('E12').Value = '  -14.64%  '
$ws.Range('D13').Value = "'0.0000248"
$ws.Range('D13').ClearFormats()
$ws.Range('E13').Value = '  -12.14%  '
$ws.Range('D14').Value = "'8.92"
$ws.Range('D14').ClearFormats()
$ws.Range('E14').Value = '  -13.93%  '
$ws.Range('D15').Value = '3.725.05'
$ws.Range('E15').Value = '  -10.69%  '
$ws.Range('E16').Value = '  -7.21%  '
$ws.Range('D17').Value = '3.214.24'
$ws.Range('E17').Value = '  -10.73%  '
$ws.Range('D18').Value = '62.577.62'
$ws.Range('E18').Value = '  -8.09%  '
$ws.Range('D19').Value = "'16.92"
$ws.Range('D19').ClearFormats()
$ws.Range('E19').Value = '  -12.03%  '
$ws.Range('D20').Value = "'10.61"
$ws.Range('D20').ClearFormats()
$ws.Range('E20').Value = '  -13.33%  '
$ws.Range('D21').Value = "'0.920"
$ws.Range('D21').ClearFormats()
$ws.Range('E21').Value = '  -13.76%  '
$ws.Range('D22').Value = "'362.41"
$ws.Range('D22').ClearFormats()
$ws.Range('E22').Value = '  -10.60%  '
$ws.Range('D23').Value = "'78.17"
$ws.Range('D23').ClearFormats()
$ws.Range('E23').Value = '  -7.92%  '
$ws.Range('D24').Value = "'10.74"
$ws.Range('D24').ClearFormats()
$ws.Range('E24').Value = '  -13.39%  '
$ws.Range('D25').Value = "'3.54"
$ws.Range('D25').ClearFormats()
$ws.Range('E25').Value = '  -15.19%  '
$ws.Range('D26').Value = "'5.90"
$ws.Range('D26').ClearFormats()
$ws.Range('E26').Value = '  -3.64%  '
$ws.Range('D27').Value = "'3.69"
$ws.Range('D27').ClearFormats()
$ws.Range('E27').Value = '  -4.25%  '
$ws.Range('D28').Value = "'2.58"
$ws.Range('D28').ClearFormats()
$ws.Range('E28').Value = '  -11.15%  '
$ws.Range('D29').Value = "'11.01"
$ws.Range('D29').ClearFormats()
$ws.Range('E29').Value = '  -12.21%  '
$ws.Range('D30').Value = "'8.06"
$ws.Range('D30').ClearFormats()
$ws.Range('E30').Value = '  -13.06%  '
$ws.Range('D31').Value = "'635.43"
$ws.Range('D31').ClearFormats()
$ws.Range('E31').Value = '  -8.72%  '
$ws.Range('D32').Value = "'27.67"
$ws.Range('D32').ClearFormats()
$ws.Range('E32').Value = '  -12.49%  '
$ws.Range('D33').Value = "'6.47"
$ws.Range('D33').ClearFormats()
$ws.Range('E33').Value = '  -15.23%  '
$ws.Range('D34').Value = "'10.88"
$ws.Range('D34').ClearFormats()
$ws.Range('E34').Value = '  -10.76%  '
$ws.Range('D35').Value = "'58.17"
$ws.Range('D35').ClearFormats()
$ws.Range('E35').Value = '  -8.62%  '
$ws.Range('D36').Value = "'1.00"
$ws.Range('D36').ClearFormats()
$ws.Range('E36').Value = '  -0.07%  '
$ws.Range('D37').Value = "'0.101"
$ws.Range('D37').ClearFormats()
$ws.Range('E37').Value = '  -11.31%  '
$ws.Range('D38').Value = "'34.99"
$ws.Range('D38').ClearFormats()
$ws.Range('E38').Value = '  -16.03%  '
$ws.Range('D39').Value = "'0.369"
$ws.Range('D39').ClearFormats()
$ws.Range('E39').Value = '  -10.51%  '
$ws.Range('D40').Value = "'0.997"
$ws.Range('D40').ClearFormats()
$ws.Range('E40').Value = '  -0.09%  '
$ws.Range('E41').Value = '  -8.64%  '
$ws.Range('D42').Value = '2.828.90'
$ws.Range('E42').Value = '  -10.83%  '
$ws.Range('D43').Value = '0.0₃0632'
$ws.Range('E43').Value = '  -16.86%  '
$ws.Range('D44').Value = "'2.60"
$ws.Range('D44').ClearFormats()
$ws.Range('E44').Value = '  -20.17%  '
$ws.Range('B45').Value = 'Fetch.AI'
$ws.Range('C45').Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range('D45').Value = "'2.32"
$ws.Range('D45').ClearFormats()
$ws.Range('E45').Value = '  -13.37%  '
$ws.Range('B46').Value = 'WEMIXToken'
$ws.Range('C46').Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range('D46').Value = "'2.55"
$ws.Range('D46').ClearFormats()
$ws.Range('E46').Value = '  -7.76%  '
$ws.Range('E47').Value = '  +0.81%  '
$ws.Range('E48').Value = '  -10.36%  '
$ws.Range('B49').Value = 'Stellar'
$ws.Range('C49').Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range('D49').Value = "'0.121"
$ws.Range('D49').ClearFormats()
$ws.Range('E49').Value = '  -8.13%  '
$ws.Range('B50').Value = 'ApeXProtocol'
$ws.Range('C50').Value = 'https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex'
$ws.Range('D50').Value = "'2.84"
$ws.Range('D50').ClearFormats()
$ws.Range('E50').Value = '  -9.29%  '
$ws.Range('D51').Value = "'129.24"
$ws.Range('D51').ClearFormats()
$ws.Range('E51').Value = '  -7.12%  '
